$d = $word.ActiveDocument

$pairs = @(
    @("31×48=1488", "12×26=312"),
    @("46×39=1794", "27×17=459"),
    @("36×57=2052", "73×31=2263"),
    @("79×18=1422", "35×33=1155"),
    @("48×16=768", "79×46=3634"),
    @("80×90=7200", "61×14=854"),
    @("52×43=2236", "51×70=3570"),
    @("85×57=4845", "52×44=2288"),
    @("31×84=2604", "41×81=3321"),
    @("92×63=5796", "16×40=640"),
    @("39×30=1170", "51×55=2805"),
    @("21×49=1029", "96×99=9504"),
    @("16×13=208", "63×40=2520"),
    @("35×17=595", "69×95=6555"),
    @("15×92=1380", "26×54=1404"),
    @("68×24=1632", "43×44=1892"),
    @("27×54=1458", "19×90=1710"),
    @("46×93=4278", "96×68=6528"),
    @("74×89=6586", "52×88=4576"),
    @("98×47=4606", "67×42=2814"),
    @("87×33=2871", "12×93=1116"),
    @("13×62=806", "16×75=1200"),
    @("81×78=6318", "57×85=4845"),
    @("30×88=2640", "11×39=429"),
    @("95×71=6745", "43×32=1376")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
